$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 4573.25
$ws.Range("I61").Value = 431
$ws.Range("J61").Value = 17000
$ws.Range("K61").Value = 1293
$ws.Range("L61").Value = 51000
$ws.Range("M61").Value = -1121
$ws.Range("N61").Value = -51344
$ws.Range("H128").Value = 109260
$ws.Range("J128").Value = 109260
$ws.Range("L128").Value = 109260
$ws.Range("N128").Value = -119220
$ws.Range("H132").Value = 2049.8857
$ws.Range("I132").Value = 1568.1818
$ws.Range("K132").Value = 4704.5454
$ws.Range("M132").Value = -2174.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2565.5
$ws.Range("I63").Value = 2478.6
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 2478.6
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -1792.6
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 2565.5
$ws.Range("I66").Value = 2478.6
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 12393
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -8961
$ws.Range("N66").Value = -21864
$ws.Range("H102").Value = 47932.375
$ws.Range("I102").Value = 51855.75
$ws.Range("J102").Value = 28315.5
$ws.Range("K102").Value = 51855.75
$ws.Range("L102").Value = 28315.5
$ws.Range("M102").Value = -50233.75
$ws.Range("N102").Value = -31559.5
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H110").Value = 2599.6
$ws.Range("I110").Value = 2166.6667
$ws.Range("J110").Value = 3249
$ws.Range("K110").Value = 2166.6667
$ws.Range("L110").Value = 3249
$ws.Range("M110").Value = -121.6667000000002
$ws.Range("N110").Value = -7339
$ws.Range("H112").Value = 11077.2
$ws.Range("J112").Value = 11077.2
$ws.Range("L112").Value = 11077.2
$ws.Range("N112").Value = -14031.2
$ws.Range("H124").Value = 46465
$ws.Range("J124").Value = 46465
$ws.Range("L124").Value = 46465
$ws.Range("N124").Value = -56285
$ws.Range("H133").Value = 78000
$ws.Range("J133").Value = 78000
$ws.Range("L133").Value = 78000
$ws.Range("N133").Value = -83060
$ws.Range("H138").Value = 69799
$ws.Range("J138").Value = 69799
$ws.Range("L138").Value = 69799
$ws.Range("N138").Value = -80079
$ws.Range("H140").Value = 62246.25
$ws.Range("J140").Value = 62246.25
$ws.Range("L140").Value = 62246.25
$ws.Range("N140").Value = -72606.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 64999
$ws.Range("J126").Value = 64999
$ws.Range("L126").Value = 64999
$ws.Range("N126").Value = -74879

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 305.1579
$ws.Range("I7").Value = 155
$ws.Range("J7").Value = 414.36365
$ws.Range("K7").Value = 155
$ws.Range("L7").Value = 414.36365
$ws.Range("M7").Value = -42
$ws.Range("N7").Value = -640.36365
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 378242.34
$ws.Range("I128").Value = 378242.34
$ws.Range("K128").Value = 1134727.02
$ws.Range("M128").Value = -1129747.02

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 1800
$ws.Range("J27").Value = 1800
$ws.Range("L27").Value = 1800
$ws.Range("N27").Value = -2132
$ws.Range("H57").Value = 21260.166
$ws.Range("I57").Value = 15833.333
$ws.Range("J57").Value = 26687
$ws.Range("K57").Value = 15833.333
$ws.Range("L57").Value = 26687
$ws.Range("M57").Value = -15013.333
$ws.Range("N57").Value = -28327
$ws.Range("H64").Value = 81499.5
$ws.Range("I64").Value = 75000
$ws.Range("K64").Value = 75000
$ws.Range("M64").Value = -74752
$ws.Range("H67").Value = 81499.5
$ws.Range("I67").Value = 75000
$ws.Range("K67").Value = 75000
$ws.Range("M67").Value = -74142
$ws.Range("H122").Value = 6252993.5
$ws.Range("I122").Value = 10002740
$ws.Range("K122").Value = 30008220
$ws.Range("M122").Value = -30005770

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1833.3334
$ws.Range("I61").Value = 1833.3334
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1833.3334
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1631.3334
$ws.Range("N61").ClearContents()
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984
$ws.Range("H93").Value = 2653.5557
$ws.Range("I93").Value = 2240
$ws.Range("J93").Value = 2771.7144
$ws.Range("K93").Value = 2240
$ws.Range("L93").Value = 2771.7144
$ws.Range("M93").Value = -992
$ws.Range("N93").Value = -5267.7144
$ws.Range("H113").Value = 1833.3334
$ws.Range("I113").Value = 1833.3334
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1833.3334
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 336.6666
$ws.Range("N113").ClearContents()
$ws.Range("H119").Value = 29684
$ws.Range("J119").Value = 29684
$ws.Range("L119").Value = 29684
$ws.Range("N119").Value = -39360
$ws.Range("H122").Value = 133337800
$ws.Range("I122").Value = 200004420
$ws.Range("K122").Value = 600013260
$ws.Range("M122").Value = -600010810

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2704
$ws.Range("I107").Value = 2049.8
$ws.Range("K107").Value = 6149.400000000001
$ws.Range("M107").Value = -4229.400000000001
$ws.Range("H109").Value = 50377
$ws.Range("J109").Value = 50377
$ws.Range("L109").Value = 50377
$ws.Range("N109").Value = -53151
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
